$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3503
$ws.Range("I62").Value = 2754.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2754.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2130.5
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3503
$ws.Range("I65").Value = 2754.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 13772.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -10652.5
$ws.Range("N65").Value = -31240

$ws.Range("H112").Value = 3981.6365
$ws.Range("J112").Value = 3979.8
$ws.Range("L112").Value = 11939.4
$ws.Range("N112").Value = -14155.4

$ws.Range("H116").Value = 5113.2856
$ws.Range("I116").Value = 4298
$ws.Range("J116").Value = 5724.75
$ws.Range("K116").Value = 4298
$ws.Range("L116").Value = 5724.75
$ws.Range("M116").Value = -856
$ws.Range("N116").Value = -12608.75

$ws.Range("H132").Value = 1810.85
$ws.Range("I132").Value = 1567.6111
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4702.8333
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2172.8333
$ws.Range("N132").Value = -17060

$ws.Range("H137").Value = 1649.9412
$ws.Range("I137").Value = 1426.9231
$ws.Range("K137").Value = 4280.7693
$ws.Range("M137").Value = -1730.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2750.549
$ws.Range("I74").Value = 2006.1892
$ws.Range("K74").Value = 2006.1892
$ws.Range("M74").Value = -1132.1892

$ws.Range("H77").Value = 2750.549
$ws.Range("I77").Value = 2006.1892
$ws.Range("K77").Value = 10030.946
$ws.Range("M77").Value = -5662.946

$ws.Range("H101").Value = 70000
$ws.Range("J101").Value = 70000
$ws.Range("L101").Value = 70000
$ws.Range("N101").Value = -76490

$ws.Range("H132").Value = 4994.5938
$ws.Range("I132").Value = 4453.72
$ws.Range("J132").Value = 6926.2856
$ws.Range("K132").Value = 13361.16
$ws.Range("L132").Value = 20778.8568
$ws.Range("M132").Value = -10831.16
$ws.Range("N132").Value = -25838.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3554.5
$ws.Range("J20").Value = 3868.3572
$ws.Range("L20").Value = 3868.3572
$ws.Range("N20").Value = -4362.3572

$ws.Range("H134").Value = 6858.846
$ws.Range("I134").Value = 7014
$ws.Range("K134").Value = 21042
$ws.Range("M134").Value = -18507

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4012.2273
$ws.Range("I31").Value = 2913.45
$ws.Range("K31").Value = 2913.45
$ws.Range("M31").Value = -2618.45

$ws.Range("H34").Value = 4012.2273
$ws.Range("I34").Value = 2913.45
$ws.Range("K34").Value = 2913.45
$ws.Range("M34").Value = -2711.45

$ws.Range("H58").Value = 10582.533
$ws.Range("I58").Value = 8033.857
$ws.Range("J58").Value = 12812.625
$ws.Range("K58").Value = 8033.857
$ws.Range("L58").Value = 12812.625
$ws.Range("M58").Value = -7830.857
$ws.Range("N58").Value = -13218.625

$ws.Range("H136").Value = 10582.533
$ws.Range("I136").Value = 8033.857
$ws.Range("J136").Value = 12812.625
$ws.Range("K136").Value = 24101.571
$ws.Range("L136").Value = 38437.875
$ws.Range("M136").Value = -21551.571
$ws.Range("N136").Value = -43537.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 217.58824
$ws.Range("J40").Value = 183
$ws.Range("L40").Value = 732
$ws.Range("N40").Value = -870

$ws.Range("H103").Value = 650
$ws.Range("J103").Value = 300
$ws.Range("L103").Value = 900
$ws.Range("N103").Value = -2658

$ws.Range("H113").Value = 727.4375
$ws.Range("I113").Value = 649.25
$ws.Range("J113").Value = 753.5
$ws.Range("K113").Value = 1947.75
$ws.Range("L113").Value = 2260.5
$ws.Range("M113").Value = 222.25
$ws.Range("N113").Value = -6600.5

$ws.Range("H122").Value = 512.9091
$ws.Range("J122").Value = 504.2
$ws.Range("L122").Value = 4537.8
$ws.Range("N122").Value = -9437.799999999999

$ws.Range("H128").Value = 518359.8
$ws.Range("I128").Value = 518359.8
$ws.Range("K128").Value = 1555079.4
$ws.Range("M128").Value = -1550099.4

$ws.Range("H140").Value = 435243.34
$ws.Range("I140").Value = 1055.9048
$ws.Range("K140").Value = 3167.7144
$ws.Range("M140").Value = 2012.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("N33").Value = -25504

$ws.Range("H44").Value = 19994
$ws.Range("I44").Value = 19994
$ws.Range("K44").Value = 19994
$ws.Range("M44").Value = -19398

$ws.Range("H70").Value = 11792.934
$ws.Range("I70").Value = 6554.8
$ws.Range("J70").Value = 14412
$ws.Range("K70").Value = 6554.8
$ws.Range("L70").Value = 14412
$ws.Range("M70").Value = -6284.8
$ws.Range("N70").Value = -14952

$ws.Range("H73").Value = 11792.934
$ws.Range("I73").Value = 6554.8
$ws.Range("J73").Value = 14412
$ws.Range("K73").Value = 6554.8
$ws.Range("L73").Value = 14412
$ws.Range("M73").Value = -5618.8
$ws.Range("N73").Value = -16284

$ws.Range("H93").Value = 38363
$ws.Range("J93").Value = 38363
$ws.Range("L93").Value = 38363
$ws.Range("N93").Value = -42107

$ws.Range("H132").Value = 2496.5
$ws.Range("I132").Value = 2497.25
$ws.Range("J132").Value = 2495
$ws.Range("K132").Value = 7491.75
$ws.Range("L132").Value = 7485
$ws.Range("M132").Value = -4961.75
$ws.Range("N132").Value = -12545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2724.75
$ws.Range("I22").Value = 2300
$ws.Range("J22").Value = 3999
$ws.Range("K22").Value = 2300
$ws.Range("L22").Value = 3999
$ws.Range("M22").Value = -2005
$ws.Range("N22").Value = -4589

$ws.Range("H27").Value = 2724.75
$ws.Range("I27").Value = 2300
$ws.Range("J27").Value = 3999
$ws.Range("K27").Value = 2300
$ws.Range("L27").Value = 3999
$ws.Range("M27").Value = -2193
$ws.Range("N27").Value = -4213

$ws.Range("H46").Value = 11174.429
$ws.Range("I46").Value = 4073.6667
$ws.Range("J46").Value = 14879.174
$ws.Range("K46").Value = 4073.6667
$ws.Range("L46").Value = 14879.174
$ws.Range("M46").Value = -3885.6667
$ws.Range("N46").Value = -15255.174

$ws.Range("H55").Value = 666.8823
$ws.Range("I55").Value = 726.1818
$ws.Range("K55").Value = 726.1818
$ws.Range("M55").Value = -553.1818

$ws.Range("H100").Value = 6768.5386
$ws.Range("I100").Value = 5624.5
$ws.Range("J100").Value = 7277
$ws.Range("K100").Value = 5624.5
$ws.Range("L100").Value = 7277
$ws.Range("M100").Value = -5083.5
$ws.Range("N100").Value = -8359

$ws.Range("H127").Value = 89000
$ws.Range("J127").Value = 89000
$ws.Range("L127").Value = 89000
$ws.Range("N127").Value = -98920

$ws.Range("H136").Value = 5744.647
$ws.Range("I136").Value = 4928.5713
$ws.Range("J136").Value = 9553
$ws.Range("K136").Value = 14785.7139
$ws.Range("L136").Value = 28659
$ws.Range("M136").Value = -12235.7139
$ws.Range("N136").Value = -33759

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5924.1787
$ws.Range("I132").Value = 5249
$ws.Range("K132").Value = 15747
$ws.Range("M132").Value = -13217
